# Update "想去人数" (F column) figures across the three data sheets to
# reflect the newly generated output (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 6896
$ws1.Range("F5").Value  = 450
$ws1.Range("F6").Value  = 154
$ws1.Range("F7").Value  = 6662
$ws1.Range("F8").Value  = 63
$ws1.Range("F10").Value = 1290
$ws1.Range("F13").Value = 400
$ws1.Range("F16").Value = 392
$ws1.Range("F17").Value = 46
$ws1.Range("F18").Value = 9
$ws1.Range("F19").Value = 5024
$ws1.Range("F20").Value = 101
$ws1.Range("F22").Value = 466
$ws1.Range("F24").Value = 199

# ---- Sheet: 演出 ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 46

# ---- Sheet: 全部类型 ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6896
$ws4.Range("F5").Value  = 0
$ws4.Range("F6").Value  = 154
$ws4.Range("F7").Value  = 6662
$ws4.Range("F8").Value  = 63
$ws4.Range("F9").Value  = 200
$ws4.Range("F13").Value = 400
$ws4.Range("F16").Value = 392
$ws4.Range("F17").Value = 46
$ws4.Range("F20").Value = 5024
$ws4.Range("F22").Value = 101
$ws4.Range("F23").Value = 120
$ws4.Range("F24").Value = 466
$ws4.Range("F25").Value = 212
$ws4.Range("F26").Value = 199
